$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 621, pushing existing rows 621-650 down to 622-651
$ws.Rows.Item(621).Insert()

# Populate the newly inserted row 621 with the new record
$ws.Cells.Item(621, 1).Value = 5
$ws.Cells.Item(621, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(621, 3).Value = "Maule"
$ws.Cells.Item(621, 4).Value = 45041
$ws.Cells.Item(621, 5).Value = 7
$ws.Cells.Item(621, 6).Value = 100112043
$ws.Cells.Item(621, 7).Value = "Pepino ensalada"
$ws.Cells.Item(621, 8).Value = "Sin especificar"
$ws.Cells.Item(621, 9).Value = "Primera"
$ws.Cells.Item(621, 10).Value = 300
$ws.Cells.Item(621, 11).Value = 12000
$ws.Cells.Item(621, 12).Value = 12000
$ws.Cells.Item(621, 13).Value = 12000
$ws.Cells.Item(621, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(621, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(621, 16).Value = 200
$ws.Cells.Item(621, 17).Value = 60
$ws.Cells.Item(621, 18).Value = "Hortaliza"
